$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.592.00'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.607.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.57%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.17'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.515'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.78%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '26.66'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.249'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.76%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0596'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0913'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.839.14'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.613.46'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.716.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.538'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.12%  '

$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.25'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.66%  '

$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '240.92'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.97%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0692'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.998'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.75%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.59'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.62%  '

$ws.Range("E27").Value = '  +4.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.36'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0471'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.21%  '

$ws.Range("E31").Value = '  +0.48%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.440.69'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.09'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.29%  '

$ws.Range("E35").Value = '  -0.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.28%  '

$ws.Range("E37").Value = '  +1.83%  '

$ws.Range("E38").Value = '  +0.04%  '

$ws.Range("E39").Value = '  +2.53%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.532'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.81%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.93'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.17%  '

$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.792'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.84%  '

$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '52.67'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +22.70%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0468'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.90%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.98%  '

$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.749.73'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.68%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.27'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.79'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.77%  '

$ws.Range("E50").Value = '  -3.44%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0103'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.39%  '
